# STS IR Bot Performer - Config.xlsx update
# "Performer. Click on OpenReturn was improved. Organized worfklows.
#  LA State Balancing second version."
#
# This inserts a new block of "LA State Balancing" constants (second
# version) just above the existing "SC State Balancing" block on the
# Constants sheet, pushing the SC block (and everything below it) down
# by 7 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# Make room: insert 7 blank rows right above the old row 120
# (StateBalancing_SC_PathCsvTemp), which then becomes row 127.
$ws.Rows("119:125").Insert()

# The inserted rows should keep the same 14.25pt custom row height used
# throughout this sheet.
$ws.Rows("119:125").RowHeight = 14.25

# --- New "LA State Balancing" (2nd version) constants -----------------
$ws.Cells.Item(120, 1).Value = "StateBalancing_LA_GrossSalesClickOnText"
$ws.Cells.Item(120, 2).Value = "GROSS SALES OF TANGIBLE###Gross Sales of Tangible###Gross sales of tangible###Gross Sales of tangible"

$ws.Cells.Item(121, 1).Value = "StateBalancing_LA_RegexGetMonthlyCheckbox"
$ws.Cells.Item(121, 2).Value = "(?<=\n.*)\w+(?=.*Monthly)"

$ws.Cells.Item(122, 1).Value = "StateBalancing_LA_ListNonLocalStates"
$ws.Cells.Item(122, 2).Value = "LA R-1029,LA R-1029E,LA R-1031"

$ws.Cells.Item(124, 1).Value = "StateBalancing_LA_TableColumns"

$ws.Cells.Item(123, 1).Value = "StateBalancing_LA_TemporalCsvFile"
$ws.Cells.Item(123, 2).Value = "Data\TemporalFileLAStateBalancing.csv"

$ws.Cells.Item(125, 1).Value = "StateBalancing_LA_ExceptionMessageQuarterlyReturn"
$ws.Cells.Item(125, 2).Value = "This return is Quarterly. Gross Sales are:"

$ws.Cells.Item(124, 2).Value = "ReturnType,FilingType,LegalEntity,ReturnName,CustomerName,GrossSales"

$ws.Cells.Item(120, 3).Value = "For LA State Balancing 2nd version"

$ws.Cells.Item(119, 1).Value = "StateBalancing_LA_GrossSalesTab"
$ws.Cells.Item(119, 2).Value = 1

# Update the selection to reflect where the author ended up after the edit.
$ws.Range("A119").Select()
